$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.413.48"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").Value = "3.554.70"
$ws.Range("E3").Value = "  +1.10%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.79"
$ws.Range("E5").Value = "  +6.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.82"
$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D8").Value = "3.551.87"
$ws.Range("E8").Value = "  +1.21%  "

$ws.Range("E10").Value = "  +4.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.25"
$ws.Range("E11").Value = "  +7.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.583"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.52"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("D15").Value = "4.136.42"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("E16").Value = "  -1.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "614.09"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "3.564.75"
$ws.Range("E18").Value = "  +1.65%  "

$ws.Range("D19").Value = "70.570.56"
$ws.Range("E19").Value = "  +2.20%  "

$ws.Range("E20").Value = "  -2.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.35"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.879"
$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.39"
$ws.Range("E23").Value = "  -15.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.70"
$ws.Range("E24").Value = "  -1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.67"
$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.46"
$ws.Range("E29").Value = "  +2.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.01"
$ws.Range("E30").Value = "  -3.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.47"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("E32").Value = "  -3.40%  "

$ws.Range("E33").Value = "  -1.78%  "

$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "571.71"
$ws.Range("E35").Value = "  -9.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.64"
$ws.Range("E36").Value = "  +6.89%  "

$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.80"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.22"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0470"
$ws.Range("E40").Value = "  +5.47%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.141"
$ws.Range("E42").Value = "  +3.56%  "

$ws.Range("D43").Value = "3.378.96"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.95"
$ws.Range("E45").Value = "  +0.39%  "

$ws.Range("E46").Value = "  +7.51%  "

$ws.Range("D47").Value = "0.0₃0700"
$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("E48").Value = "  +2.11%  "

$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.81"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("E51").Value = "  -0.02%  "
